$p = $ppt.ActivePresentation

# Slide 1: subtitle placeholder - merge "Bài " + "23. " + "Thu " into a single run
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange.Characters(1, 12)
$tr1.Text = "Bài 23. Thu "

# Slide 28: title placeholder - split "Bài tập" into "Bài " + "tập 23.1"
$s28 = $p.Slides.Item(28)
$shp28 = $s28.Shapes.Item(2)
$tr28 = $shp28.TextFrame.TextRange.Characters(5, 3)
$tr28.Text = "tập 23.1"
